# Auto-generated Excel COM-interop edit script
# Applies the cell-value changes described by the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3742549
$ws.Range("J40").Value = 5166144.5
$ws.Range("L40").Value = 5166144.5
$ws.Range("N40").Value = -5166494.5
$ws.Range("H64").Value = 90582910
$ws.Range("J64").Value = 4732.3335
$ws.Range("L64").Value = 4732.3335
$ws.Range("N64").Value = -5228.3335
$ws.Range("H67").Value = 90582910
$ws.Range("J67").Value = 4732.3335
$ws.Range("L67").Value = 4732.3335
$ws.Range("N67").Value = -6448.3335
$ws.Range("H74").Value = 68185680
$ws.Range("I74").Value = 375002240
$ws.Range("J74").Value = 4222.222
$ws.Range("K74").Value = 375002240
$ws.Range("L74").Value = 4222.222
$ws.Range("M74").Value = -375001304
$ws.Range("N74").Value = -6094.222
$ws.Range("H77").Value = 68185680
$ws.Range("I77").Value = 375002240
$ws.Range("J77").Value = 4222.222
$ws.Range("K77").Value = 1875011200
$ws.Range("L77").Value = 21111.11
$ws.Range("M77").Value = -1875006520
$ws.Range("N77").Value = -30471.11
$ws.Range("H132").Value = 158929.84
$ws.Range("I132").Value = 230406.53
$ws.Range("K132").Value = 691219.59
$ws.Range("M132").Value = -688689.59
$ws.Range("H138").Value = 2554.3953
$ws.Range("I138").Value = 1963.8422
$ws.Range("J138").Value = 2721.8657
$ws.Range("K138").Value = 5891.5266
$ws.Range("L138").Value = 8165.597099999999
$ws.Range("M138").Value = -751.5266000000001
$ws.Range("N138").Value = -18445.5971
$ws.Range("H140").Value = 72332.375
$ws.Range("J140").Value = 72564.28999999999
$ws.Range("L140").Value = 72564.28999999999
$ws.Range("N140").Value = -82924.28999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1789
$ws.Range("I74").Value = 1239.8572
$ws.Range("K74").Value = 1239.8572
$ws.Range("M74").Value = -365.8571999999999
$ws.Range("H77").Value = 1789
$ws.Range("I77").Value = 1239.8572
$ws.Range("K77").Value = 6199.286
$ws.Range("M77").Value = -1831.286
$ws.Range("H97").Value = 746.53845
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 458131.1
$ws.Range("I102").Value = 490622.6
$ws.Range("K102").Value = 490622.6
$ws.Range("M102").Value = -489000.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 767
$ws.Range("I22").Value = 786.375
$ws.Range("J22").Value = 689.5
$ws.Range("K22").Value = 786.375
$ws.Range("L22").Value = 689.5
$ws.Range("M22").Value = -613.375
$ws.Range("N22").Value = -1035.5
$ws.Range("H82").Value = 34814.25
$ws.Range("I82").Value = 16419
$ws.Range("J82").Value = 90000
$ws.Range("K82").Value = 16419
$ws.Range("L82").Value = 90000
$ws.Range("M82").Value = -16036
$ws.Range("N82").Value = -90766
$ws.Range("H85").Value = 34814.25
$ws.Range("I85").Value = 16419
$ws.Range("J85").Value = 90000
$ws.Range("K85").Value = 16419
$ws.Range("L85").Value = 90000
$ws.Range("M85").Value = -15093
$ws.Range("N85").Value = -92652
$ws.Range("H116").Value = 72300
$ws.Range("J116").Value = 72300
$ws.Range("L116").Value = 72300
$ws.Range("N116").Value = -81478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 304.9
$ws.Range("I7").Value = 229.8
$ws.Range("J7").Value = 380
$ws.Range("K7").Value = 229.8
$ws.Range("L7").Value = 380
$ws.Range("M7").Value = -116.8
$ws.Range("N7").Value = -606
$ws.Range("H29").Value = 800
$ws.Range("J29").Value = 800
$ws.Range("L29").Value = 800
$ws.Range("N29").Value = -1386
$ws.Range("H31").Value = 5167.2104
$ws.Range("I31").Value = 1750.8889
$ws.Range("J31").Value = 6743.974
$ws.Range("K31").Value = 1750.8889
$ws.Range("L31").Value = 6743.974
$ws.Range("M31").Value = -1455.8889
$ws.Range("N31").Value = -7333.974
$ws.Range("H34").Value = 5167.2104
$ws.Range("I34").Value = 1750.8889
$ws.Range("J34").Value = 6743.974
$ws.Range("K34").Value = 1750.8889
$ws.Range("L34").Value = 6743.974
$ws.Range("M34").Value = -1548.8889
$ws.Range("N34").Value = -7147.974
$ws.Range("H58").Value = 2596.9395
$ws.Range("I58").Value = 1419.96
$ws.Range("K58").Value = 1419.96
$ws.Range("M58").Value = -1216.96
$ws.Range("H136").Value = 2596.9395
$ws.Range("I136").Value = 1419.96
$ws.Range("K136").Value = 4259.88
$ws.Range("M136").Value = -1709.88
$ws.Range("H141").Value = 80730.75
$ws.Range("J141").Value = 92657.16
$ws.Range("L141").Value = 92657.16
$ws.Range("N141").Value = -103017.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 509171.62
$ws.Range("J33").Value = 2002052
$ws.Range("L33").Value = 12012312
$ws.Range("N33").Value = -12012878
$ws.Range("H68").Value = 1298963
$ws.Range("J68").Value = 1298963
$ws.Range("L68").Value = 3896889
$ws.Range("N68").Value = -3898511
$ws.Range("H71").Value = 1298963
$ws.Range("J71").Value = 1298963
$ws.Range("L71").Value = 11690667
$ws.Range("N71").Value = -11698779
$ws.Range("H86").Value = 500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 1500
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -3872
$ws.Range("H89").Value = 500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 4500
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -16356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 393.44446
$ws.Range("I2").Value = 220.57895
$ws.Range("J2").Value = 804
$ws.Range("K2").Value = 220.57895
$ws.Range("L2").Value = 804
$ws.Range("M2").Value = -107.57895
$ws.Range("N2").Value = -1030
$ws.Range("H11").Value = 25025000
$ws.Range("I11").Value = 26700000
$ws.Range("K11").Value = 26700000
$ws.Range("M11").Value = -26699861
$ws.Range("H107").Value = 7936997.5
$ws.Range("I107").Value = 7936997.5
$ws.Range("K107").Value = 7936997.5
$ws.Range("M107").Value = -7935077.5
$ws.Range("H140").Value = 29999.5
$ws.Range("J140").Value = 89998
$ws.Range("L140").Value = 89998
$ws.Range("N140").Value = -100358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6500.067
$ws.Range("I46").Value = 4375
$ws.Range("J46").Value = 6827
$ws.Range("K46").Value = 4375
$ws.Range("L46").Value = 6827
$ws.Range("M46").Value = -4187
$ws.Range("N46").Value = -7203
$ws.Range("H68").Value = 2068554.5
$ws.Range("I68").Value = 3788746
$ws.Range("J68").Value = 4324.6
$ws.Range("K68").Value = 3788746
$ws.Range("L68").Value = 4324.6
$ws.Range("M68").Value = -3787997
$ws.Range("N68").Value = -5822.6
$ws.Range("H71").Value = 2068554.5
$ws.Range("I71").Value = 3788746
$ws.Range("J71").Value = 4324.6
$ws.Range("K71").Value = 18943730
$ws.Range("L71").Value = 21623
$ws.Range("M71").Value = -18939986
$ws.Range("N71").Value = -29111
$ws.Range("H82").Value = 3473269.5
$ws.Range("I82").Value = 7813047
$ws.Range("J82").Value = 1447.8
$ws.Range("K82").Value = 7813047
$ws.Range("L82").Value = 1447.8
$ws.Range("M82").Value = -7812686
$ws.Range("N82").Value = -2169.8
$ws.Range("H85").Value = 3473269.5
$ws.Range("I85").Value = 7813047
$ws.Range("J85").Value = 1447.8
$ws.Range("K85").Value = 7813047
$ws.Range("L85").Value = 1447.8
$ws.Range("M85").Value = -7811799
$ws.Range("N85").Value = -3943.8
$ws.Range("H122").Value = 10559.571
$ws.Range("J122").Value = 17823.834
$ws.Range("L122").Value = 53471.50199999999
$ws.Range("N122").Value = -58371.50199999999
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4237.9165
$ws.Range("I126").Value = 3805
$ws.Range("K126").Value = 11415
$ws.Range("M126").Value = -8945
$ws.Range("H132").Value = 71430780
$ws.Range("I132").Value = 1750
$ws.Range("K132").Value = 5250
$ws.Range("M132").Value = -2720
$ws.Range("H136").Value = 8628.692999999999
$ws.Range("I136").Value = 3597.5173
$ws.Range("K136").Value = 10792.5519
$ws.Range("M136").Value = -8242.5519
